$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 734.5
$ws.Range("I2").Value = 481.5
$ws.Range("K2").Value = 481.5
$ws.Range("M2").Value = -368.5
$ws.Range("H17").Value = 3191.1333
$ws.Range("J17").Value = 3191.1333
$ws.Range("L17").Value = 9573.3999
$ws.Range("N17").Value = -9909.3999
$ws.Range("H33").Value = 488.05264
$ws.Range("I33").Value = 174.25
$ws.Range("K33").Value = 174.25
$ws.Range("M33").Value = 54.75
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()  # was -3075
$ws.Range("H43").Value = 2198.25
$ws.Range("I43").Value = 2198.25
$ws.Range("K43").Value = 2198.25
$ws.Range("M43").Value = -2129.25
$ws.Range("H80").Value = 543.0968
$ws.Range("I80").Value = 382.85715
$ws.Range("J80").Value = 675.05884
$ws.Range("K80").Value = 1148.57145
$ws.Range("L80").Value = 2025.17652
$ws.Range("M80").Value = -150.5714499999999
$ws.Range("N80").Value = -4021.17652
$ws.Range("H83").Value = 543.0968
$ws.Range("I83").Value = 382.85715
$ws.Range("J83").Value = 675.05884
$ws.Range("K83").Value = 3445.71435
$ws.Range("L83").Value = 6075.52956
$ws.Range("M83").Value = 1546.28565
$ws.Range("N83").Value = -16059.52956
$ws.Range("H88").Value = 2923.85
$ws.Range("J88").Value = 3126.2354
$ws.Range("L88").Value = 3126.2354
$ws.Range("N88").Value = -3938.2354
$ws.Range("H91").Value = 2923.85
$ws.Range("J91").Value = 3126.2354
$ws.Range("L91").Value = 3126.2354
$ws.Range("N91").Value = -5934.2354
$ws.Range("H112").Value = 2888.7017
$ws.Range("J112").Value = 2904.5894
$ws.Range("L112").Value = 8713.768199999999
$ws.Range("N112").Value = -10929.7682
$ws.Range("H116").Value = 5580.357
$ws.Range("I116").Value = 6995
$ws.Range("K116").Value = 6995
$ws.Range("M116").Value = -3553
$ws.Range("H132").Value = 2735.647
$ws.Range("I132").Value = 2394
$ws.Range("K132").Value = 7182
$ws.Range("M132").Value = -4652
$ws.Range("H138").Value = 5104.0386
$ws.Range("J138").Value = 6411
$ws.Range("L138").Value = 19233
$ws.Range("N138").Value = -29513

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6512.2856
$ws.Range("I61").Value = 4419.7334
$ws.Range("K61").Value = 4419.7334
$ws.Range("M61").Value = -4207.7334
$ws.Range("H74").Value = 1756.9302
$ws.Range("I74").Value = 1763.7028
$ws.Range("J74").Value = 1715.1666
$ws.Range("K74").Value = 1763.7028
$ws.Range("L74").Value = 1715.1666
$ws.Range("M74").Value = -889.7028
$ws.Range("N74").Value = -3463.1666
$ws.Range("H77").Value = 1756.9302
$ws.Range("I77").Value = 1763.7028
$ws.Range("J77").Value = 1715.1666
$ws.Range("K77").Value = 8818.513999999999
$ws.Range("L77").Value = 8575.833000000001
$ws.Range("M77").Value = -4450.513999999999
$ws.Range("N77").Value = -17311.833
$ws.Range("H110").Value = 1485.8334
$ws.Range("I110").Value = 1183.0667
$ws.Range("K110").Value = 1183.0667
$ws.Range("M110").Value = 861.9332999999999
$ws.Range("H114").Value = 86663.336
$ws.Range("J114").Value = 86663.336
$ws.Range("L114").Value = 86663.336
$ws.Range("N114").Value = -95341.336
$ws.Range("H122").Value = 5869.696
$ws.Range("I122").Value = 4999.5
$ws.Range("K122").Value = 14998.5
$ws.Range("M122").Value = -12548.5
$ws.Range("H136").Value = 6512.2856
$ws.Range("I136").Value = 4419.7334
$ws.Range("K136").Value = 13259.2002
$ws.Range("M136").Value = -10709.2002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 25026426
$ws.Range("I94").Value = 31250884
$ws.Range("K94").Value = 31250884
$ws.Range("M94").Value = -31250433
$ws.Range("H99").Value = 3601.65
$ws.Range("I99").Value = 3910.7646
$ws.Range("K99").Value = 3910.7646
$ws.Range("M99").Value = -2412.7646
$ws.Range("H132").Value = 90780
$ws.Range("J132").Value = 90780
$ws.Range("L132").Value = 90780
$ws.Range("N132").Value = -100900

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22100.4
$ws.Range("I16").Value = 20999.875
$ws.Range("J16").Value = 23358.143
$ws.Range("K16").Value = 20999.875
$ws.Range("L16").Value = 23358.143
$ws.Range("M16").Value = -20712.875
$ws.Range("N16").Value = -23932.143
$ws.Range("H94").Value = 1216
$ws.Range("I94").Value = 879.2
$ws.Range("K94").Value = 879.2
$ws.Range("M94").Value = -428.2
$ws.Range("H105").Value = 3688.0667
$ws.Range("I105").Value = 2523
$ws.Range("K105").Value = 2523
$ws.Range("M105").Value = -776
$ws.Range("H113").Value = 22100.4
$ws.Range("I113").Value = 20999.875
$ws.Range("J113").Value = 23358.143
$ws.Range("K113").Value = 20999.875
$ws.Range("L113").Value = 23358.143
$ws.Range("M113").Value = -18829.875
$ws.Range("N113").Value = -27698.143
$ws.Range("H132").Value = 1536.4348
$ws.Range("I132").Value = 1479.4147
$ws.Range("K132").Value = 4438.2441
$ws.Range("M132").Value = -1908.2441
$ws.Range("H141").Value = 294224.88
$ws.Range("J141").Value = 294224.88
$ws.Range("L141").Value = 294224.88
$ws.Range("N141").Value = -304584.88

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1017.7059
$ws.Range("I2").Value = 1683.3334
$ws.Range("K2").Value = 10100.0004
$ws.Range("M2").Value = -9987.000400000001
$ws.Range("H32").Value = 9162.5
$ws.Range("J32").Value = 9331
$ws.Range("L32").Value = 27993
$ws.Range("N32").Value = -28559
$ws.Range("H46").Value = 31566.031
$ws.Range("I46").Value = 959.05554
$ws.Range("J46").Value = 68294.39999999999
$ws.Range("K46").Value = 2877.16662
$ws.Range("L46").Value = 204883.2
$ws.Range("M46").Value = -2786.16662
$ws.Range("N46").Value = -205065.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1509.742
$ws.Range("I102").Value = 1288.963
$ws.Range("K102").Value = 1288.963
$ws.Range("M102").Value = 333.037
$ws.Range("H122").Value = 3458.6538
$ws.Range("I122").Value = 2703.8125
$ws.Range("J122").Value = 4666.4
$ws.Range("K122").Value = 8111.4375
$ws.Range("L122").Value = 13999.2
$ws.Range("M122").Value = -5661.4375
$ws.Range("N122").Value = -18899.2
$ws.Range("H132").Value = 4677.087
$ws.Range("I132").Value = 4944.476
$ws.Range("J132").Value = 1869.5
$ws.Range("K132").Value = 14833.428
$ws.Range("L132").Value = 5608.5
$ws.Range("M132").Value = -12303.428
$ws.Range("N132").Value = -10668.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 15999
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 15999
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 15999
$ws.Range("M2").ClearContents()  # was -11888
$ws.Range("N2").Value = -16223
$ws.Range("H46").Value = 2042.7142
$ws.Range("I46").Value = 1349.5
$ws.Range("J46").Value = 2320
$ws.Range("K46").Value = 1349.5
$ws.Range("L46").Value = 2320
$ws.Range("M46").Value = -1161.5
$ws.Range("N46").Value = -2696
$ws.Range("H55").Value = 300.5
$ws.Range("I55").Value = 300.5
$ws.Range("K55").Value = 300.5
$ws.Range("M55").Value = -127.5
$ws.Range("H68").Value = 4008.2307
$ws.Range("I68").Value = 3853.4443
$ws.Range("J68").Value = 4356.5
$ws.Range("K68").Value = 3853.4443
$ws.Range("L68").Value = 4356.5
$ws.Range("M68").Value = -3104.4443
$ws.Range("N68").Value = -5854.5
$ws.Range("H71").Value = 4008.2307
$ws.Range("I71").Value = 3853.4443
$ws.Range("J71").Value = 4356.5
$ws.Range("K71").Value = 19267.2215
$ws.Range("L71").Value = 21782.5
$ws.Range("M71").Value = -15523.2215
$ws.Range("N71").Value = -29270.5
$ws.Range("H122").Value = 9983.143
$ws.Range("I122").Value = 9983.143
$ws.Range("K122").Value = 29949.429
$ws.Range("M122").Value = -27499.429
$ws.Range("H136").Value = 4613.7334
$ws.Range("I136").Value = 3953.1
$ws.Range("K136").Value = 11859.3
$ws.Range("M136").Value = -9309.299999999999
$ws.Range("H140").Value = 134500
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()  # was -129770

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26666.666
$ws.Range("I54").Value = 26666.666
$ws.Range("K54").Value = 26666.666
$ws.Range("M54").Value = -26146.666
$ws.Range("H86").Value = 69991
$ws.Range("J86").Value = 69991
$ws.Range("L86").Value = 69991
$ws.Range("N86").Value = -72237
$ws.Range("H89").Value = 69991
$ws.Range("J89").Value = 69991
$ws.Range("L89").Value = 349955
$ws.Range("N89").Value = -361187
$ws.Range("H122").Value = 4363
$ws.Range("I122").Value = 4041.8333
$ws.Range("J122").Value = 4844.75
$ws.Range("K122").Value = 12125.4999
$ws.Range("L122").Value = 14534.25
$ws.Range("M122").Value = -9675.499899999999
$ws.Range("N122").Value = -19434.25
$ws.Range("H132").Value = 4470.7803
$ws.Range("I132").Value = 2437.4062
$ws.Range("K132").Value = 7312.2186
$ws.Range("M132").Value = -4782.2186
$ws.Range("H136").Value = 3217.342
$ws.Range("I136").Value = 3295.7585
$ws.Range("J136").Value = 2964.6667
$ws.Range("K136").Value = 9887.2755
$ws.Range("L136").Value = 8894.000100000001
$ws.Range("M136").Value = -7337.2755
$ws.Range("N136").Value = -13994.0001
